$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -17.88608949611059
$ws.Range("C2").Value = 1.491207895056935
$ws.Range("D2").Value = -17.88608949611059
$ws.Range("E2").Value = -17.88608949611059
$ws.Range("F2").Value = -17.88608949611059
$ws.Range("G2").Value = -17.88608949611059
$ws.Range("H2").Value = -17.88608949611059
$ws.Range("I2").Value = -17.88608949611059
$ws.Range("J2").Value = -17.88608949611059
$ws.Range("K2").Value = -17.88608949611059

$ws.Range("B3").Value = -17.88608949611059
$ws.Range("C3").Value = -17.88608949611059
$ws.Range("D3").Value = -17.88608949611059
$ws.Range("E3").Value = -17.88608949611059
$ws.Range("F3").Value = -17.88608949611059
$ws.Range("G3").Value = -17.88608949611059
$ws.Range("H3").Value = -17.88608949611059
$ws.Range("I3").Value = 2.518841034924626
$ws.Range("J3").Value = -17.88608949611059
$ws.Range("K3").Value = -17.88608949611059

$ws.Range("B4").Value = -17.88608949611059
$ws.Range("C4").Value = 2.038074950236763
$ws.Range("D4").Value = -17.88608949611059
$ws.Range("E4").Value = -17.88608949611059
$ws.Range("F4").Value = 3.991469122308193
$ws.Range("G4").Value = -17.88608949611059
$ws.Range("H4").Value = 1.761321498379276
$ws.Range("I4").Value = -17.88608949611059
$ws.Range("J4").Value = 2.034474167299327
$ws.Range("K4").Value = -17.88608949611059

$ws.Range("B5").Value = -17.88608949611059
$ws.Range("C5").Value = 1.526862755104298
$ws.Range("D5").Value = -17.88608949611059
$ws.Range("E5").Value = -17.88608949611059
$ws.Range("F5").Value = -17.88608949611059
$ws.Range("G5").Value = 2.804755091077614
$ws.Range("H5").Value = -17.88608949611059
$ws.Range("I5").Value = -17.88608949611059
$ws.Range("J5").Value = -17.88608949611059
$ws.Range("K5").Value = -17.88608949611059

$ws.Range("B6").Value = -17.88608949611059
$ws.Range("C6").Value = -17.88608949611059
$ws.Range("D6").Value = -17.88608949611059
$ws.Range("E6").Value = -17.88608949611059
$ws.Range("F6").Value = -17.88608949611059
$ws.Range("G6").Value = -17.88608949611059
$ws.Range("H6").Value = -17.88608949611059
$ws.Range("I6").Value = -17.88608949611059
$ws.Range("J6").Value = -17.88608949611059
$ws.Range("K6").Value = -17.88608949611059

$ws.Range("B7").Value = 2.563511573284007
$ws.Range("C7").Value = -17.88608949611059
$ws.Range("D7").Value = -17.88608949611059
$ws.Range("E7").Value = -17.88608949611059
$ws.Range("F7").Value = -17.88608949611059
$ws.Range("G7").Value = -17.88608949611059
$ws.Range("H7").Value = -17.88608949611059
$ws.Range("I7").Value = -17.88608949611059
$ws.Range("J7").Value = -17.88608949611059
$ws.Range("K7").Value = -17.88608949611059

$ws.Range("B8").Value = -17.88608949611059
$ws.Range("C8").Value = -17.88608949611059
$ws.Range("D8").Value = -17.88608949611059
$ws.Range("E8").Value = 2.365722270951849
$ws.Range("F8").Value = -17.88608949611059
$ws.Range("G8").Value = -17.88608949611059
$ws.Range("H8").Value = -17.88608949611059
$ws.Range("I8").Value = -17.88608949611059
$ws.Range("J8").Value = -17.88608949611059
$ws.Range("K8").Value = -17.88608949611059

$ws.Range("B9").Value = 3.816443821320298
$ws.Range("C9").Value = -17.88608949611059
$ws.Range("D9").Value = -17.88608949611059
$ws.Range("E9").Value = -17.88608949611059
$ws.Range("F9").Value = -17.88608949611059
$ws.Range("G9").Value = -17.88608949611059
$ws.Range("H9").Value = -17.88608949611059
$ws.Range("I9").Value = -17.88608949611059
$ws.Range("J9").Value = -17.88608949611059
$ws.Range("K9").Value = -17.88608949611059

$ws.Range("B10").Value = -17.88608949611059
$ws.Range("C10").Value = -17.88608949611059
$ws.Range("D10").Value = -17.88608949611059
$ws.Range("E10").Value = -17.88608949611059
$ws.Range("F10").Value = -17.88608949611059
$ws.Range("G10").Value = -17.88608949611059
$ws.Range("H10").Value = -17.88608949611059
$ws.Range("I10").Value = 0.9189647249045724
$ws.Range("J10").Value = -17.88608949611059
$ws.Range("K10").Value = 2.048046227708304

$ws.Range("B11").Value = -17.88608949611059
$ws.Range("C11").Value = -17.88608949611059
$ws.Range("D11").Value = -17.88608949611059
$ws.Range("E11").Value = 2.589589517745853
$ws.Range("F11").Value = -17.88608949611059
$ws.Range("G11").Value = 3.092166390843081
$ws.Range("H11").Value = -17.88608949611059
$ws.Range("I11").Value = -17.88608949611059
$ws.Range("J11").Value = -17.88608949611059
$ws.Range("K11").Value = 2.123270328821777

$ws.Range("B12").Value = -17.88608949611059
$ws.Range("C12").Value = -17.88608949611059
$ws.Range("D12").Value = -17.88608949611059
$ws.Range("E12").Value = -17.88608949611059
$ws.Range("F12").Value = -17.88608949611059
$ws.Range("G12").Value = -17.88608949611059
$ws.Range("H12").Value = -17.88608949611059
$ws.Range("I12").Value = -17.88608949611059
$ws.Range("J12").Value = -17.88608949611059
$ws.Range("K12").Value = -17.88608949611059

$ws.Range("B13").Value = -17.88608949611059
$ws.Range("C13").Value = -17.88608949611059
$ws.Range("D13").Value = -17.88608949611059
$ws.Range("E13").Value = 2.714706350325744
$ws.Range("F13").Value = -17.88608949611059
$ws.Range("G13").Value = -17.88608949611059
$ws.Range("H13").Value = -17.88608949611059
$ws.Range("I13").Value = -17.88608949611059
$ws.Range("J13").Value = 1.672004184964471
$ws.Range("K13").Value = 1.950461742020148

$ws.Range("B14").Value = -17.88608949611059
$ws.Range("C14").Value = -17.88608949611059
$ws.Range("D14").Value = -17.88608949611059
$ws.Range("E14").Value = -17.88608949611059
$ws.Range("F14").Value = -17.88608949611059
$ws.Range("G14").Value = -17.88608949611059
$ws.Range("H14").Value = -17.88608949611059
$ws.Range("I14").Value = -17.88608949611059
$ws.Range("J14").Value = -17.88608949611059
$ws.Range("K14").Value = 2.116666471581801

$ws.Range("B15").Value = -17.88608949611059
$ws.Range("C15").Value = -17.88608949611059
$ws.Range("D15").Value = -17.88608949611059
$ws.Range("E15").Value = -17.88608949611059
$ws.Range("F15").Value = -17.88608949611059
$ws.Range("G15").Value = -17.88608949611059
$ws.Range("H15").Value = -17.88608949611059
$ws.Range("I15").Value = -17.88608949611059
$ws.Range("J15").Value = -17.88608949611059
$ws.Range("K15").Value = -17.88608949611059

$ws.Range("B16").Value = -17.88608949611059
$ws.Range("C16").Value = -17.88608949611059
$ws.Range("D16").Value = -17.88608949611059
$ws.Range("E16").Value = -17.88608949611059
$ws.Range("F16").Value = -17.88608949611059
$ws.Range("G16").Value = -17.88608949611059
$ws.Range("H16").Value = -17.88608949611059
$ws.Range("I16").Value = -17.88608949611059
$ws.Range("J16").Value = 1.801991737290511
$ws.Range("K16").Value = -17.88608949611059

$ws.Range("B17").Value = -17.88608949611059
$ws.Range("C17").Value = 2.986541050662487
$ws.Range("D17").Value = 4.321922437061814
$ws.Range("E17").Value = -17.88608949611059
$ws.Range("F17").Value = -17.88608949611059
$ws.Range("G17").Value = -17.88608949611059
$ws.Range("H17").Value = 1.561213724964474
$ws.Range("I17").Value = 2.405406120221441
$ws.Range("J17").Value = 2.531669296249595
$ws.Range("K17").Value = -17.88608949611059

$ws.Range("B18").Value = -17.88608949611059
$ws.Range("C18").Value = -17.88608949611059
$ws.Range("D18").Value = -17.88608949611059
$ws.Range("E18").Value = -17.88608949611059
$ws.Range("F18").Value = -17.88608949611059
$ws.Range("G18").Value = -17.88608949611059
$ws.Range("H18").Value = 2.038538913447388
$ws.Range("I18").Value = 1.584197996945708
$ws.Range("J18").Value = 1.785432566551379
$ws.Range("K18").Value = -17.88608949611059

$ws.Range("B19").Value = -17.88608949611059
$ws.Range("C19").Value = -17.88608949611059
$ws.Range("D19").Value = -17.88608949611059
$ws.Range("E19").Value = -17.88608949611059
$ws.Range("F19").Value = -17.88608949611059
$ws.Range("G19").Value = -17.88608949611059
$ws.Range("H19").Value = 1.394438889735923
$ws.Range("I19").Value = 1.331605565106415
$ws.Range("J19").Value = -17.88608949611059
$ws.Range("K19").Value = -17.88608949611059

$ws.Range("B20").Value = -17.88608949611059
$ws.Range("C20").Value = 0.07721866829048495
$ws.Range("D20").Value = -17.88608949611059
$ws.Range("E20").Value = -17.88608949611059
$ws.Range("F20").Value = 2.033601658299822
$ws.Range("G20").Value = -17.88608949611059
$ws.Range("H20").Value = 1.705880190120029
$ws.Range("I20").Value = 0.646057127019674
$ws.Range("J20").Value = -17.88608949611059
$ws.Range("K20").Value = 1.724994754066239

$ws.Range("B21").Value = -17.88608949611059
$ws.Range("C21").Value = 0.2862256882901275
$ws.Range("D21").Value = -17.88608949611059
$ws.Range("E21").Value = 1.177571536494001
$ws.Range("F21").Value = -17.88608949611059
$ws.Range("G21").Value = 2.16503578586761
$ws.Range("H21").Value = 1.871608139236682
$ws.Range("I21").Value = -17.88608949611059
$ws.Range("J21").Value = -17.88608949611059
$ws.Range("K21").Value = -17.88608949611059

